$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.4
$ws.Range("A3").Value = 90
$ws.Range("B3").Value = 0
$ws.Range("A4").Value = 120
$ws.Range("B4").Value = 0

$ws.Range("C7").Select()
